# Applies the "fixed all list and new tag associations" style cleanup:
#   - List1 / List3change / List4change: turn off contextual spacing
#   - List6 / List7 / List8: drop the (incorrect) Heading3 basis, restore the
#     proper run formatting (explicit fonts, bold/bCs flags and complex-script
#     size) that the List6Char/List7Char/List8Char sibling styles already had
#   - List7 / List7Char: bold is actually turned ON for style 7 (the others
#     stay non-bold)

$d = $word.ActiveDocument
$styles = $d.Styles

# ---- List1 : remove contextual spacing -----------------------------------
$styles.Item("List1").NoSpaceBetweenParagraphsOfSameStyle = $false

# ---- List6 -----------------------------------------------------------------
$s6 = $styles.Item("List6")
$s6.BaseStyle = ""
$s6.NoSpaceBetweenParagraphsOfSameStyle = $false
$s6.Font.NameAscii = "Times New Roman"
$s6.Font.NameOther = "Times New Roman"
$s6.Font.NameBi = "Times New Roman"
$s6.Font.BoldBi = $true
$s6.Font.SizeBi = 12

# ---- List7 -------------------------------------------------------------
$s7 = $styles.Item("List7")
$s7.BaseStyle = ""
$s7.Font.Bold = $true
$s7.Font.NameAscii = "Times New Roman"
$s7.Font.NameOther = "Times New Roman"
$s7.Font.NameBi = "Times New Roman"
$s7.Font.SizeBi = 12

# ---- List7Char ---------------------------------------------------------
$s7c = $styles.Item("List7Char")
$s7c.Font.Bold = $true
$s7c.Font.BoldBi = $false

# ---- List8 -------------------------------------------------------------
$s8 = $styles.Item("List8")
$s8.BaseStyle = ""
$s8.NoSpaceBetweenParagraphsOfSameStyle = $false
$s8.Font.NameAscii = "Times New Roman"
$s8.Font.NameOther = "Times New Roman"
$s8.Font.NameBi = "Times New Roman"
$s8.Font.BoldBi = $true
$s8.Font.SizeBi = 12

# ---- List3change / List4change : remove contextual spacing ----------------
$styles.Item("List3change").NoSpaceBetweenParagraphsOfSameStyle = $false
$styles.Item("List4change").NoSpaceBetweenParagraphsOfSameStyle = $false
